# Simulated Wild Card round and logged it
# Update the 2021 Target Depth Data for OFF (offense) and DEF (defense) sheets
# with the additional game's stats folded into the season totals.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")

# Row 2 (H)
$wsOff.Range("B2").Value = 252
$wsOff.Range("C2").Value = 186
$wsOff.Range("D2").Value = 75
$wsOff.Range("E2").Value = 39

# Row 3 (R)
$wsOff.Range("B3").Value = 275
$wsOff.Range("C3").Value = 191
$wsOff.Range("D3").Value = 63
$wsOff.Range("E3").Value = 30
$wsOff.Range("F3").Value = 4

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")

# Row 2 (H)
$wsDef.Range("B2").Value = 279
$wsDef.Range("C2").Value = 206
$wsDef.Range("D2").Value = 79
$wsDef.Range("E2").Value = 27
$wsDef.Range("G2").Value = 4

# Row 3 (R)
$wsDef.Range("B3").Value = 234
$wsDef.Range("C3").Value = 171
$wsDef.Range("D3").Value = 44
$wsDef.Range("E3").Value = 19

$wb.Save()
